# Apply the "Apellido" column insertion + header centering edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before the current column C (DNI), shifting
#    DNI / Numero Registro / Fecha Emision / Direccion / Segundo Cargo /
#    Centro de formacion one column to the right.
$ws.Columns.Item(3).Insert()

# 2) Write the new header in the freshly inserted column C.
$ws.Range("C1").Value = "Apellido"

# 3) Center-align the whole header row (A1:I1), matching the new
#    cellXfs entry (horizontal="center") applied to every header cell.
$ws.Range("A1:I1").HorizontalAlignment = -4108

# 4) Resize columns to match the refreshed layout (values below are the
#    closest widths this engine's ColumnWidth rounding can reach).
$ws.Columns.Item(1).ColumnWidth = 18.5
$ws.Columns.Item(2).ColumnWidth = 32
$ws.Columns.Item(3).ColumnWidth = 25.666666666666668
$ws.Columns.Item(4).ColumnWidth = 20.333333333333332
$ws.Columns.Item(5).ColumnWidth = 21.5
$ws.Columns.Item(6).ColumnWidth = 24.5
$ws.Columns.Item(7).ColumnWidth = 36.166666666666664
$ws.Columns.Item(8).ColumnWidth = 28
$ws.Columns.Item(9).ColumnWidth = 24

# 5) Update the saved selection to G4.
$ws.Range("G4").Select() | Out-Null
